$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B ("data") cells for existing rows 2-6: convert from numeric
#     date serials into plain literal text "15/08/2024" (no special
#     number-format any more).
$ws.Range("B2:B6").ClearFormats()
$ws.Range("B2").Value = "15/08/2024"
$ws.Range("B3").Value = "15/08/2024"
$ws.Range("B4").Value = "15/08/2024"
$ws.Range("B5").Value = "15/08/2024"
$ws.Range("B6").Value = "15/08/2024"

# --- Swap the product descriptions between rows 3 and 4
$ws.Range("E3").Value = "SACOS PARA HOT DOG"
$ws.Range("E4").Value = "SACOS PARA HAMBURGÃO"

# --- Append the new daily entries (rows 7-18)
$newRows = @(
    @(5,  "16/08/2024", 2068354, "LUVAS PLÁSTICAS",               "BOMPACK",  "100Un.", 13),
    @(6,  "16/08/2024", 2068354, "PAPEL MANTEIGA",                 "BOMPACK",  "4m",     13),
    @(7,  "16/08/2024", 2068354, "PAPEL ALUMINIO",                 "BOMPACK",  "4m",     13),
    @(8,  "16/08/2024", 2068358, "SACOS PLÁSTICOS",                "BOT",      "20Un.",  13),
    @(9,  "16/08/2024", 2068358, "SACOS PARA LIXO",                "BOM",      "5Un.",   13),
    @(10, "16/08/2024", 2068358, "SACOS PLÁSTICOS P/ ALIMENTOS",   "BOT",      "50Un.",  13),
    @(11, "16/08/2024", 2068361, "SACOS PARA LIXO",                "BOM",      "50Un.",  13),
    @(12, "16/08/2024", 2068361, "SACOS PARA LIXO",                "BOT",      "10Un.",  13),
    @(13, "16/08/2024", 2068361, "SACOS PARA LIXO",                "BOT",      "10Un.",  13),
    @(14, "16/08/2024", 2068362, "SACOS PARA LIXO",                "BOT",      "4Un.",   13),
    @(15, "16/08/2024", 2068363, "SACOS PARA LIXO",                "FAM ROLL", "30Un.",  13),
    @(16, "16/08/2024", 2068363, "SACOS PARA LIXO",                "FAM ROLL", "30Un.",  13)
)

$r = 7
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    # A lone apostrophe stores as an empty (but still text-typed) cell,
    # matching the blank "tipo_saida" placeholder already used in D2:D6.
    $ws.Cells.Item($r, 4).Value = "'"
    $ws.Cells.Item($r, 4).ClearFormats()
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# The A column (index) cells reuse the bordered/centered style already used
# by A2:A6 - copy it across (format-only paste) instead of re-deriving font
# / border / alignment by hand, so no *new* style entries get created.
$ws.Range("A2").Copy()
$ws.Range("A7:A18").PasteSpecial(-4122)
